$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite header row: proper capitalization, no leading space
$ws.Range("B1").Value = "Pld"
$ws.Range("C1").Value = "Won"
$ws.Range("D1").Value = "Lost"
$ws.Range("E1").Value = "Tied"
$ws.Range("F1").Value = "Net RR"
$ws.Range("G1").Value = "Pts"

# Fix the team name in row 9 (was incorrectly duplicated "Delhi Capitals")
$ws.Range("A9").Value = "Kings XI Punjab"

# Select the table range, as reflected in the saved view state
$ws.Range("A1:G9").Select()
